# Flip the sign of the material-recycled figures (B2:E5) on every yearly
# sheet: every non-zero numeric value in that block becomes negative
# (positive -> negative), matching the commit's "+ to -" sign change.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 5; $row++) {
        for ($col = 2; $col -le 5; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $val = $cell.Value()
            if ($val -ne 0) {
                $cell.Value = -$val
            }
        }
    }
}
